$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:F1) ---
# C1 "Marca" -> "Concentração_Obtida"
$ws.Range("C1").Value = "Concentração_Obtida"
# Insert new "Laboratório" header in D1, shift old D/E headers right into E/F
$ws.Range("D1").Value = "Laboratório"
$ws.Range("E1").Value = "Registro"
$ws.Range("F1").Value = "PDF"
# New F1 header cell should carry the same header formatting as the rest
# of row 1 (bold, centered, bordered) -- copy it from the neighbouring cell.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows: add Concentração_Obtida (C), shift Laboratório into D, ---
# --- keep Registro in E, and set the new PDF/status column (F) ---

# Row 2 (Item 6)
$ws.Range("C2").Value = "1000ml"
$ws.Range("D2").Value = "JP INDUSTRIA FARMACEUTICA S/A"
$ws.Range("F2").Value = "Pendente"

# Row 3 (Item 7)
$ws.Range("C3").Value = "250ml"
$ws.Range("D3").Value = "JP INDUSTRIA FARMACEUTICA S/A"
$ws.Range("F3").Value = "Pendente"

# Row 4 (Item 8)
$ws.Range("C4").Value = "500ml"
$ws.Range("D4").Value = "JP INDUSTRIA FARMACEUTICA S/A"
$ws.Range("F4").Value = "Pendente"

# Row 5 (Item 9)
$ws.Range("C5").Value = "1000ml"
$ws.Range("D5").Value = "JP INDUSTRIA FARMACEUTICA S/A"
$ws.Range("F5").Value = "Pendente"

# --- Registro column (E2:E5) keeps its original text values ("104910019" / ---
# --- "104910020"), but these are numeric-looking strings, so the range is ---
# --- formatted as text first or Excel would silently store them as numbers. ---
$regRange = $ws.Range("E2:E5")
$regRange.NumberFormat = "@"
$ws.Range("E2").Value = "104910019"
$ws.Range("E3").Value = "104910019"
$ws.Range("E4").Value = "104910019"
$ws.Range("E5").Value = "104910020"
# Reset back to the workbook's default (unformatted) style now that the
# values are locked in as text, so no visible formatting change remains.
$regRange.Style = "Normal"
